$d = $word.ActiveDocument

# 1. Merge the split runs of the "El módulo Informes..." paragraph into a
#    single run by replacing the whole (unchanged) text via Find & Replace.
$oldText = "El módulo Informes es el módulo principal de la aplicación. Contribuye para generar respuesta a las inquietudes del Administrador respecto a la cantidad de insumos que tiene en la tienda, la cantidad de transacciones hechas y el balance generado."
$newText = "El módulo Informes es el módulo principal de la aplicación. Contribuye para generar respuesta a las inquietudes del Administrador respecto a la cantidad de insumos que tiene en la tienda, la cantidad de transacciones hechas y el balance generado."
$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)

# 2. Remove the _GoBack bookmark from its old location first, while it is
#    still the only one in the document, so a later by-name lookup can't
#    resolve to the wrong (freshly inserted) instance.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 3. Locate that paragraph again (its contents are now a single run) and
#    insert a new paragraph right after it.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "El módulo Informes es el módulo principal*") {
        $target = $p
        break
    }
}
$target.Range.InsertParagraphAfter()
$newPara = $target.Next()

# 4. Inject the new "Módulo Ventas" heading paragraph plus a following
#    empty paragraph that holds the relocated _GoBack bookmark.
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Ttulo2"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">Módulo </w:t></w:r><w:r><w:t>Ventas</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:bookmarkStart w:id="7" w:name="_GoBack"/><w:bookmarkEnd w:id="7"/></w:p>'
$newPara.Range.InsertXML($xml)

# 5. Update the header relationship + create footnotes/endnotes parts:
#    add a footnote and an endnote, then remove the references again. This
#    mirrors the placeholder footnotes.xml/endnotes.xml parts left behind.
$r0 = $d.Range(0, 0)
$fn = $d.Footnotes.Add($r0, "", "x")
$fn.Reference.Delete()
$en = $d.Endnotes.Add($r0, "", "x")
$en.Reference.Delete()
